# Add newly-extracted otolith records (rows 216-231) to Sheet1.
# Columns: A=Trawl, B=Species, C=Length, D=No. Extracted, E=Well No., F=Plate No.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(216, 4, "Maurolicus muelleri",    "47mm SL", 2, "G7",  "TC2"),
    @(217, 4, "Maurolicus muelleri",    "46mm SL", 2, "G8",  "TC2"),
    @(218, 4, "Maurolicus muelleri",    "49mm SL", 2, "G9",  "TC2"),
    @(219, 4, "Maurolicus muelleri",    "48mm SL", 2, "G10", "TC2"),
    @(220, 4, "Maurolicus muelleri",    "50mm SL", 2, "G11", "TC2"),
    @(221, 4, "Maurolicus muelleri",    "52mm SL", 2, "G12", "TC2"),
    @(222, 4, "Maurolicus muelleri",    "54mm SL", 2, "H1",  "TC2"),
    @(223, 6, "Electrona risso",        "22mm SL", 2, "H10", "TC2"),
    @(224, 6, "Electrona risso",        "21mm SL", 2, "H11", "TC2"),
    @(225, 6, "Electrona risso",        "NR",      2, "H12", "TC2"),
    @(226, 6, "Argyropelecus olfersii", "43mm SL", 1, "A1",  "TC4"),
    @(227, 6, "Argyropelecus olfersii", "47mm SL", 1, "A2",  "TC4"),
    @(228, 6, "Argyropelecus olfersii", "49mm SL", 2, "A3",  "TC4"),
    @(229, 6, "Argyropelecus olfersii", "43mm SL", 2, "A4",  "TC4"),
    @(230, 6, "Argyropelecus olfersii", "52mm SL", 1, "A5",  "TC4"),
    @(231, 6, "Argyropelecus olfersii", "NR",      2, "A6",  "TC4")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 2).Font.Italic = $true
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

$ws.Range("D231").Select() | Out-Null
